$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for 2025-08-23 (serial 45892), which had no data recorded.
# This shifts all subsequent rows up by one, matching the new date sequence.
$ws.Rows(3).Delete()

# A couple of the (now shifted-up) rows get corrected evening readings.
$ws.Range("C3").Value = 106.75
$ws.Range("E3").Value = 28.3
$ws.Range("C5").Value = 106.05
$ws.Range("E5").Value = 28.9

# Clear out the medication-dose column for all data rows - dosing is no
# longer being logged.
$ws.Range("F2:F7").ClearContents()
